# Workbook was re-saved after correcting the product name (missing a dash
# after "199") and leaving the workbook positioned on the ProductLoanOutput
# sheet/cell, matching the target commit.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Correct the product name value (adds the missing hyphen) on both sheets.
$correctedName = "199-MS-EI-DB-DL-REC-NON-RNI-CTRFD-SAR-MD-TR-1-ONTIME"
$wsInput.Range("B1").Value  = $correctedName
$wsOutput.Range("B1").Value = $correctedName

# Reposition the selection on the input sheet away from row 6 to B1.
[void]$wsInput.Range("B1").Select()

# Make ProductLoanOutput the active sheet/tab, with B1 selected there too.
[void]$wsOutput.Activate()
[void]$wsOutput.Range("B1").Select()
